$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'Langit Penuh Harapan'
$ws.Cells.Item(2, 3).Value = 'Di cahaya pagi, kita akan menemukan jalan kita Melalui malam yang gelap, hari yang lebih cerah Kita akan berpegang pada harapan, dan tidak pernah menyerah Dan hati kita akan berdetak, dengan cinta yang dalam Kita akan bangkit di atas, kebisingan dan kesakitan Dan menemukan suara kita, dalam kesunyian lagi Kita akan berdiri bersama, tangan di tangan Dan bersinar seperti bintang, di negeri yang penuh harapan Mimpi kita akan terbang, dengan sayap elang Dan semangat kita akan, bernyanyi selamanya Kita akan menjadi perubahan, yang kita ingin lihat Dan membawa cinta dan damai, kepada kamu dan saya Kita akan menjadi cahaya, di malam yang gelap Dan memimpin satu sama lain, ke cahaya yang lebih cerah Kita tidak akan pernah menyerah, kita tidak akan pernah kehilangan iman Dan langit penuh harapan kita, akan tetap selamanya'
$ws.Cells.Item(2, 4).Value = 'semua usia'

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'Halaman yang Belum Ditulis'
$ws.Cells.Item(3, 3).Value = 'Di buku kehidupan, kita menulis cerita kita sendiri Dengan setiap langkah, bab baru terbuka Kita adalah penulis, dari takdir kita sendiri Dan pilihan kita akan, membentuk sejarah kita Kita akan membalik halaman, pada setiap hari baru Dan mengisi ruang kosong, dengan cara kita sendiri Kita akan menulis tentang cinta, tentang harapan, tentang mimpi Dan mengisi halaman yang belum ditulis, dengan tema yang paling liar Kita akan mengambil pena, dan menulis nama kita sendiri Dan menandatangani karya kita, dengan hati yang tidak terkalahkan Kita akan mengisi halaman, dengan setiap warna yang cerah Dan membuat kesan kita, pada kanvas dunia malam ini Kita akan menulis tentang keberanian, tentang hati, tentang jiwa Dan mengisi halaman yang belum ditulis, dengan cerita yang harus diceritakan Kita akan membuat jalan kita sendiri, melalui yang tidak diketahui Dan menulis bestseller, yang akan selamanya ditampilkan'
$ws.Cells.Item(3, 4).Value = 'semua usia'

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'Sinar-Sinar Kesunyian'
$ws.Cells.Item(4, 3).Value = 'Berjalan di pantai, merasa begitu bebas Matahari bersinar terang, hanya kita berdua Kita mengejar burung camar, tertawa dengan gembira Membuat kenangan, hanya kita berdua Hidup adalah perjalanan, kita baru memulai Setiap langkah kita ambil, hati kita terus berputar Kita akan menari dalam hujan, kita akan bernyanyi dalam matahari Ikatan kita akan tumbuh, cinta kita akan menang Kita akan membentangkan sayap, kita akan terbang tinggi Kita akan mencapai ketinggian baru, kita tidak akan pernah menyerah Kita akan bersinar seperti bintang, pada malam yang cerah Kita akan berkilau seperti berlian, dalam cahaya pagi Kita akan menjadi perubahan, kita akan membuatnya nyata Kita akan menjadi harapan, kita akan membuatnya terasa Kita akan menjadi cinta, yang akan selamanya ada'
$ws.Cells.Item(4, 4).Value = 'semua usia'

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Warna-Warna Kehidupan'
$ws.Cells.Item(5, 3).Value = 'Merah adalah mawar, yang mekar di hatiku Kuning adalah matahari, yang bersinar sejak awal Biru adalah langit, yang kulihat setiap hari Hijau adalah rumput, yang kulalui di jalanku Jingga adalah senja, yang kukejar setiap malam Ungu adalah mimpi, yang kugenggam erat Hidup adalah kanvas, dengan warna yang begitu cerah Kita adalah seniman, yang melukis cahaya kita sendiri Kita akan mencampur dan mencocokkan, kita akan membuat warna kita sendiri Kita akan membuat jalanku sendiri, kita akan melihatnya sampai akhir Kita akan menambahkan kilau, kita akan menambahkan gaya Kita akan membuat karya agung, yang tak tertandingi Kita akan menjadi pelangi, setelah malam yang berbadai Kita akan bersinar begitu terang, dengan semua kekuatan kita Kita akan menjadi warna, yang membuat hidup begitu cerah'
$ws.Cells.Item(5, 4).Value = 'semua usia'

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'Serenade Cahaya Bintang'
$ws.Cells.Item(6, 3).Value = 'Pada malam hari, ketika bintang muncul Saya akan bernyanyi, dengan hati yang jernih Melodi itu, beresonansi dekat Simfoni, yang membawa kita bergembira Dengan setiap nada, hati kita berpadu Dalam harmoni sempurna, kita akan bersinar Seperti pertunjukan celestial, begitu ilahi Cinta kita akan berkilau, sepanjang waktu Dalam kesunyian, cinta kita akan tumbuh Dengan setiap detak, hati kita akan bercahaya Dalam cahaya bintang, cinta kita akan bersinar Selamanya dan selalu, cinta kita akan menjadi milikku Pada malam hari, ketika bintang muncul Cinta kita akan berkilau, tahun demi tahun Dengan setiap nada, cinta kita akan tumbuh Dalam cahaya bintang, hati kita akan bercahaya'
$ws.Cells.Item(6, 4).Value = 'semua usia'

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Keajaiban Dunia'
$ws.Cells.Item(7, 3).Value = 'Dari gunung tinggi, ke lautan lebar Ada keindahan di mana-mana, kita tidak dapat bersembunyi Keajaiban dunia, tidak pernah berhenti Harta karun, yang tidak pernah habis Dari gurun panas, ke hutan dalam Dunia penuh, dengan rahasia yang kita simpan Keajaiban dunia, tidak pernah memudar Karpet ajaib, yang belum kita tukar Dengan setiap langkah, cerita baru terbuka Perjalanan penemuan, yang tidak pernah tua Keajaiban dunia, tidak pernah tidur Mimpi yang menjadi kenyataan, yang kita dapat simpan Dari gunung tinggi, ke lautan lebar Keajaiban dunia, akan menjadi petunjuk kita Dengan setiap langkah, cerita baru terbuka Keajaiban dunia, tidak pernah tua'
$ws.Cells.Item(7, 4).Value = 'semua usia'

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'Pemimpi Kecil'
$ws.Cells.Item(8, 3).Value = 'Di negeri berfantasi, aku bebas berkelana Di mana imajinasi tumbuh, dan mimpi dibuat di rumah Aku seorang pemimpi kecil, dengan hati yang cerah Aku akan mengejar bintang malam ini, dan membuat segalanya baik-baik saja Aku akan terbang dengan sayap elang, dan menyentuh embun pagi Aku akan menari di ujung pelangi, dan menyanyikan lagu atau dua Di negeri berfantasi, aku bebas menjadi Seorang pemimpi kecil, liar dan bebas Aku akan membangun istana tinggi, dengan menara menyentuh langit Aku akan berlayar di tujuh lautan, dan menemukan harta karun tersembunyi mengapa Aku akan menjadi seorang eksplorator berani, dengan hati yang sejati Aku akan membuat teman-teman baru di sepanjang jalan, dan belajar melihatnya melalui Di negeri berfantasi, aku bebas bermain Aku akan tertawa dan mencintai dan bermimpi, setiap hari Aku akan menjadi seorang pemimpi kecil, dengan hati yang cerah Berkilau seperti bintang, pada malam yang ajaib'
$ws.Cells.Item(8, 4).Value = 'anak'

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Taman Persahabatan'
$ws.Cells.Item(9, 3).Value = 'Di taman persahabatan, kita akan menanam benih Di mana cinta dan kebaikan tumbuh, dan hati dibebaskan Kita akan menyiraminya dengan tawa, dan sinar matahari setiap hari Dan menonton ikatan kita tumbuh lebih kuat, dengan cara yang istimewa Kita akan menambahkan warna-warna cerah, dengan bunga-bunga dari setiap warna Dan membuat karya agung, yang unik dan sejati Di taman persahabatan, kita akan membuat kenangan juga Dan menghargai setiap momen, yang dibagi antara aku dan kamu Kita akan membuat janji kepada satu sama lain, untuk selalu ada di sana Melalui tebal dan tipis, melalui sukacita dan pedih Di taman persahabatan, kita akan tumbuh dan berkembang Bersama-sama selamanya, berdampingan Kita akan membuat teman-teman baru, dan menambahkan ke kru kita Dan mengisi taman dengan cinta, yang kuat dan sejati Di taman persahabatan, kita akan berkembang dan mekar juga Dengan ikatan yang tidak terputus, yang selalu baru'
$ws.Cells.Item(9, 4).Value = 'anak'

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'Negara Imajinasi'
$ws.Cells.Item(10, 3).Value = 'Di dunia fantasi dan mimpi Di mana imajinasi tidak tahu skema Kita bisa menjadi siapa saja yang kita mau Di negeri ini, kita bebas Kita bisa terbang, kita bisa menyelam Kita bisa menjelajah, dan meminta lebih Imajinasi kita adalah kunci Untuk membuka rahasia kita dan kamu Di negeri ini, kita bebas berkelana Kita bisa menjadi apa saja yang kita lihat Dari naga ke unicorn, dan lebih Imajinasi kita adalah pintu Ke dunia kemungkinan yang tak terbatas Di mana mimpi menjadi kenyataan, dan sihir adalah realita Di Negara Imajinasi, kita bebas berkelana Di mana kreativitas tidak tahu rumah'
$ws.Cells.Item(10, 4).Value = 'anak'

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Pengelana Kecil'
$ws.Cells.Item(11, 3).Value = 'Saya pengelana kecil, dalam petualangan besar Saya menemukan hal-hal baru, dengan setiap langkah dan setiap ukuran Saya belajar tentang dunia, dan semua keajaibannya Saya menemukan teman-teman baru, dan membuat kenangan yang abadi Saya mendaki gunung tinggi, dan menyeberangi sungai lebar Saya menjelajah yang tidak diketahui, dengan hati penuh kebanggaan Saya pengelana kecil, dalam misi besar Saya menemukan jalanku, dan melakukan yang terbaik Saya membuka rahasia, dan memecahkan misteri juga Saya pengelana kecil, dengan hati yang benar-benar biru Saya meninggalkan jejak, dan meninggalkan bekas Saya pengelana kecil, di ruang yang luas ini Saya belajar dan tumbuh, dengan setiap penemuan baru Saya pengelana kecil, satu-satunya Saya berani dan saya gagah, dan saya tidak pernah takut Saya pengelana kecil, dalam petualangan setiap hari'
$ws.Cells.Item(11, 4).Value = 'anak'

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Mata Berkilau'
$ws.Cells.Item(12, 3).Value = 'Di cahaya pagi, aku melihat wajahmu Senyummu begitu cerah, tempat yang bahagia Tertawamu bergema, suara yang menyenangkan Harta yang ditemukan, hati yang murni dan bulat Kamu menari dan berputar, dengan langkah yang bebas Jiwa yang bebas, hati yang penuh dengan kegembiraan Mata berkilaumu, bersinar seperti matahari Harta yang berharga, yang harus dimenangkan Di mimpi-mimpimu, semoga hal-hal yang indah terjadi Semoga hatimu, dipenuhi dengan cinta untuk dipegang Semoga jalammu, diterangi dengan cahaya yang membimbing Semoga semangatmu, terbang malam ini Kamu adalah bintang, yang bersinar sangat terang Sebuah mercu suar, di kegelapan malam Kamu adalah hadiah, harta yang sangat langka Jiwa yang berharga, yang tidak ada bandingannya'
$ws.Cells.Item(12, 4).Value = 'anak'

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'Rapsodi Pelangi'
$ws.Cells.Item(13, 3).Value = 'Di langit, pelangi muncul Pemandangan yang berwarna, yang membawa kita bersuka Merah, oranye, kuning, hijau dan biru Simfoni, warna yang baru Keajaiban pelangi, membawa kita pergi Ke negeri yang ajaib, di mana mimpi-mimpi datang bermain Kita menari dan bernyanyi, dengan hati yang gembira Dalam irama pelangi, kita menemukan bagian kita Dengan setiap warna, sebuah cerita diungkapkan Tentang harapan dan cinta, yang tidak pernah tua Keindahan pelangi, adalah harta yang langka Simbol janji, yang selalu ada Jadi mari kita kejar pelangi, dan ikuti petunjuknya Ke dunia yang ajaib, di mana mimpi-mimpi dibebaskan Dan dalam warnanya, semoga hati kita menjadi ringan'
$ws.Cells.Item(13, 4).Value = 'anak'

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'Gema di Malam'
$ws.Cells.Item(14, 3).Value = 'Dalam kesunyian, saya mendengar suara Anda Rahasia yang di bisikkan, pilihan yang tulus Bintang di atas, berkelip-kelip terang Pertunjukan celestial, pada malam yang tak berujung Dunia adalah gelap, tapi saya tidak sendiri Bayangan menari, pada irama jiwa saya Saya akan mengikuti irama, ke yang tidak diketahui Di mana gema hati saya, akan ditampilkan Dalam kesunyian, saya akan menemukan jalan Melalui kegelapan, fajar baru akan bergoyang Malam akan hilang, dan matahari akan terbit Hari baru, dengan mata yang terbuka Gema di malam, akan menunjukkan saya jalan Ke tempat di mana cinta, akan bersinar selamanya Dan dalam kesunyian, saya akan menemukan suara saya Melodi harapan, yang akan membuat beberapa kebisingan'
$ws.Cells.Item(14, 4).Value = 'remaja'

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'Jiwa Bunga Liar'
$ws.Cells.Item(15, 3).Value = 'Di taman mimpi, saya akan menemukan tempat saya Di mana bunga liar, mekar dengan kemuliaan Kelopaknya lembut, warnanya cerah Keindahan yang langka, di pagi hari Saya akan membiarkan semangat saya, liar dan bebas Seperti bunga liar, yang dimaksudkan untuk bebas Saya akan menari dalam angin, dengan hati yang terbuka Dan membiarkan jiwa saya, menjadi karya seni Di taman mimpi, saya akan menemukan suara saya Melodi harapan, yang akan membuat beberapa kebisingan Saya akan membiarkan jiwa bunga liar saya, bersinar terang Sebuah mercu suar, di kegelapan malam Saya akan menyebar akar, dalam bumi Dan membiarkan semangat saya, memberikan kelahiran baru Untuk dunia yang menunggu, untuk awal yang baru Di mana cinta dan harapan, akan selalu di hati kita'
$ws.Cells.Item(15, 4).Value = 'remaja'

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'Langit Tengah Malam'
$ws.Cells.Item(16, 3).Value = 'Dalam kesunyian malam Saya mendengar bintang di atas Bisikan cahaya jauh Membimbing saya ke sisi lain Di mana mimpi dibuat dari api Dan detak jantung seperti drum Saya akan mengejar bayangan, saya akan berhasrat Menemukan kebenaran, untuk memiliki kesenangan Di langit tengah malam, saya akan menemukan jalan Melalui kegelapan, saya akan menyinari cahaya Saya akan bangkit, saya akan merebut hari Dan semangat saya akan terbang Saya akan menari di bawah laut berbintang Dengan setiap langkah, saya akan bebas Di langit tengah malam, saya akan menemukan diri'
$ws.Cells.Item(16, 4).Value = 'remaja'

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'Hati yang Tak Terhenti'
$ws.Cells.Item(17, 3).Value = 'Saya memiliki api yang membakar dalam Sebuah nyala yang berkedip, tapi tidak pernah mati Itu adalah detak hati yang tak terhenti saya Yang mendorong saya maju, tepat dari awal Saya akan bangkit, saya akan naik ke atas Saya akan membiarkan suara saya terdengar, saya akan membalik halaman Saya akan mengejar mimpi, saya akan membuatnya nyata Saya tidak akan pernah menyerah, saya tidak akan pernah kalah Hati saya berdegup, itu keras dan jelas Saya akan mengikuti passion saya, saya akan mengusir rasa takut Saya akan berlari sejauh, saya akan menyeberangi garis Saya akan membuat tanda, saya akan bersinar Saya akan tak terhenti, saya akan ilahi Dengan setiap langkah, saya akan membuatnya milik saya'
$ws.Cells.Item(17, 4).Value = 'remaja'

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'Mimpi Elektrik'
$ws.Cells.Item(18, 3).Value = 'Dalam kesunyian, saya mendengar jantung saya berdetak keras Irama yang bergema, api yang membakar dengan bangga Saya menutup mata, dan saya melihat cahaya Dunia kemungkinan, di mana saya terbang Saya mengejar tinggi, saya merasa hidup Mimpi elektrik, itulah di mana saya tumbuh Saya tidak terhenti, saya sedang naik Semangat saya meningkat, saya membuka mata saya Saya bebas dari rantai, yang mengikat saya Saya membuka sayap, saya meninggalkan masa lalu Saya bersinar terang, seperti bintang di malam Saya tidak terhenti, saya merasa cahaya Saya elektrik, saya hidup, saya berapi Saya tidak terhenti, hati saya adalah keinginan saya Saya mengejar mimpi, saya merasa panas Saya elektrik, saya hidup, saya lengkap'
$ws.Cells.Item(18, 4).Value = 'remaja'

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'Hilang dalam Momen'
$ws.Cells.Item(19, 3).Value = 'Saya berdiri di tepi, dari hari yang baru Dunia penuh dengan keajaiban, dalam setiap cara Saya merasa hidup, saya merasa bebas Saya hilang dalam momen, di mana saya seharusnya Saya melepaskan, semua ketakutan saya Saya menghadapi yang tidak diketahui, melalui semua air mata saya Saya menemukan suara saya, saya berbicara pikiran saya Saya mematahkan rantai, yang mengikat saya kembali dalam waktu Saya bersinar terang, seperti mercu suar di malam Saya tidak terhenti, saya merasa cahaya Saya hilang dalam momen, di mana cinta adalah petunjuk Saya menemukan jalan saya, ke sisi lain Saya melepaskan, semua keraguan saya Saya menghadapi keindahan, yang semua tentang Saya hilang dalam momen, di mana saya bebas untuk menjadi Saya tidak terhenti, itulah saya'
$ws.Cells.Item(19, 4).Value = 'remaja'

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'Bisikan dalam Kegelapan'
$ws.Cells.Item(20, 3).Value = 'Di bayangan, saya mendengar suara Anda Sebuah bisikan lembut, sebuah pilihan yang tulus Untuk mengikuti jalan, yang belum diketahui Dan menemukan keberanian, untuk membuatnya milik saya Kegelapan menyelimuti, tapi saya tidak takut Karena dalam kesunyian, saya mendengar hati saya berkata Bahwa saya kuat, dan saya berani Dan saya akan menemukan jalannya, melalui gelombang kegelapan Bisikan dalam kegelapan, mereka membimbing saya Dan membawa saya ke cahaya, yang bersinar kembali Saya akan memegang harapan, dan tidak akan melepaskan Karena dalam kegelapan, saya akan menemukan cara untuk tumbuh Dan meskipun jalan itu panjang, dan berliku-liku Saya akan mengikuti bisikan, yang berbicara kebenaran kepada Anda Dan saya akan bangkit di atas, bayangan yang jatuh Dan saya akan berdiri teguh, dan menghadapi semuanya'
$ws.Cells.Item(20, 4).Value = 'dewasa'

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 'Api yang Padam'
$ws.Cells.Item(21, 3).Value = 'Api yang dulu berkobar, kini hanya sebuah percikan Sebuah api yang padam, dari cinta yang telah hilang Nyala api yang dulu menari, kini hanya sebuah kilau Sebuah api yang mati, yang telah kehilangan cara untuk tumbuh Saya ditinggalkan dengan abu, dari apa yang dulu ada Sebuah kenangan cinta, yang telah kehilangan melodi Kehangatan yang dulu menghibur, kini hanya sebuah kedinginan Sebuah api yang padam, yang telah kehilangan kehendak Saya akan memegang masa lalu, dan menghargai rasa sakit Karena dalam api yang padam, saya akan menemukan cinta yang tetap Meskipun api itu padam, cinta itu akan tetap ada Sebuah api yang membakar, yang tidak akan pernah hilang Dan meskipun nyala api itu mati, cinta itu tidak akan pernah padam Karena dalam api yang padam, cinta kita akan tetap dibuat Dan saya akan memegang itu, dan tidak akan melepaskan Karena dalam api yang padam, cinta kita akan tetap bersinar'
$ws.Cells.Item(21, 4).Value = 'dewasa'

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'Gelombang Sunyi'
$ws.Cells.Item(22, 3).Value = 'Dalam kesunyian malam Saya mendengar bisikan masa lalu Angin lembut yang berbisik rendah Tentang kenangan yang akan tetap ada Bintang di atas, mereka berkelipan terang Pertunjukan celestial, pemandangan yang indah Dunia ini sunyi, tempat yang damai Di mana mimpi dan kenyataan berkelindan Dalam jam ini, saya menemukan kedamaian Perasaan tenang, dunia yang berhenti Gelombang sunyi, mereka menghantam pantai Suara yang menenangkan, saya minta lebih Dalam kesunyian ini, saya bebas berkelana Di mana hati saya dapat sembuh, jiwa saya dapat pulang Dan dalam kesunyian, saya dibuat baru'
$ws.Cells.Item(22, 4).Value = 'dewasa'

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'Selamanya Ditemukan'
$ws.Cells.Item(23, 3).Value = 'Dalam kedalaman jiwa saya, cahaya menyala Sebuah mercu suar di kegelapan, hati yang milik saya Suara yang berbisik rendah, ''kamu cukup'' Cinta yang mengangkat saya, dan memberi saya kekuatan Dunia luar, itu memudar pergi Saya menemukan jalan saya, melalui hari yang paling gelap Saya mencari jawaban, di bintang di atas Dan menemukan tujuan saya, dalam desah kesepian Dalam kesunyian, saya menemukan hati saya Cinta yang murni, karya seni Selamanya ditemukan, dalam kedalaman saya Cinta yang nyata, cinta yang bebas Cinta yang menyembuhkan, cinta yang memperbaiki Cinta yang milik saya, hingga akhir waktu Dan dalam kehangatan itu, saya dibuat lengkap'
$ws.Cells.Item(23, 4).Value = 'dewasa'

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 'Bayang-Bayang Kenangan'
$ws.Cells.Item(24, 3).Value = 'Di senja pikiran saya Terletak dunia waktu yang terlupakan Di mana bayang-bayang menari dan desah berbisik Dan hantu cinta berlalu Di dalam kesunyian, saya mendengar suara Anda Sebuah melodi yang menghantui, pilihan yang tulus Untuk memegang apa yang kita temukan Atau membiarkannya hilang, seperti butiran pasir Di bayang-bayang, saya akan menemukan jalan Melalui kegelapan, ke hari yang baru Di mana cahaya cinta akan bersinar sangat terang Dan membimbing saya melalui malam yang paling gelap Dan meskipun kenangan mungkin memudar Di dalam hati saya, cinta kita tidak akan pernah beranjak Karena di bayang-bayang kenangan, saya akan menemukan jalan Kembali kepada Anda, kembali ke irama cinta yang manis'
$ws.Cells.Item(24, 4).Value = 'dewasa'

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 'Horizon yang Berkobar'
$ws.Cells.Item(25, 3).Value = 'Di tepi fajar, saya berdiri Di mana dunia terbangun, berpegangan tangan Dengan matahari terbit, hati saya berdegup kencang Sebagai horizon memanggil, saya mengambil masa lalu Di dalam api jiwa saya, mimpi terbang tinggi Sebuah passion yang membakar, yang membimbing saya melalui malam Melalui api ketakutan, saya akan menemukan jalan Ke fajar yang baru, ke hari yang baru Di mana dunia penuh dengan kemungkinan yang tak terbatas Dan horizon membentang, liar dan bebas Saya akan mengejar matahari, saya akan mengejar mimpi Dan membiarkan semangat saya terbang tinggi, seperti teriakan elang Karena di horizon yang berkobar, saya akan menemukan rumah Di mana api hati saya, akan selamanya berkeliaran Dan di dalam api, saya akan menemukan jalan Ke diri saya yang baru, ke hari yang baru'
$ws.Cells.Item(25, 4).Value = 'dewasa'
